$d = $word.ActiveDocument

# Locate the word "teit" in the document without relying on hard-coded offsets.
$found = $d.Content
$ok = $found.Find.Execute("teit", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($ok) {
    $start = $found.Start
    $end   = $found.End

    # 1) Append ".." right after "teit" (this becomes its own run below).
    $rTail = $d.Range($end, $end)
    $rTail.InsertBefore("..")

    # 2) Turn the leading "t" into "T", forcing it into its own run by
    #    toggling a character property on and back off (Word keeps the
    #    resulting run boundary even though the formatting ends up
    #    unchanged, since adjacent runs with identical rPr would
    #    otherwise be silently re-merged).
    $rFirst = $d.Range($start, $start + 1)
    $rFirst.Text = "T"
    $rFirst.Bold = 1
    $rFirst.Bold = 0

    # 3) Split the newly-appended ".." into its own run the same way.
    $rDots = $d.Range($end, $end + 2)
    $rDots.Bold = 1
    $rDots.Bold = 0
}
